# 133_GESTPROJ_G01.xlsx — "Add files via upload"
#
# Journal de travail now records a new day (06.04.2023, serial 45022) of
# work covering "théorie" (0.5h) and "REST 1 réglement de problème" (2.5h);
# the 31.03.2023 (C15) entry is corrected from 1h45 to 2h30. The weekly
# total (C22) recalculates automatically from the SUM formula already in
# the sheet.
#
# On the Planning sheet, the "ServiceRest2" column (H) is marked with an
# "X" in the Théorie + TT row (row 22), same as the ServiceRest1 column
# (G) right next to it.
#
# Finally the last on-screen selection is nudged: Planning ends up
# scrolled/selected at G24, and Journal de travail ends up selected at
# B19 (and remains the active tab, as it already was).

$wb = $excel.ActiveWorkbook

$planning = $wb.Worksheets.Item("Planning")
$journal  = $wb.Worksheets.Item("Journal de travail")

# --- Planning: flag ServiceRest2 (H22) the same way ServiceRest1 (G22) is flagged ---
$planning.Range("H22").Value = "X"

# --- Journal de travail: fix the 31.03.2023 entry and log the new day ---
$journal.Range("C15").Value = 2.5

$journal.Range("A16").Value = 45022
$journal.Range("B16").Value = "théorie"
$journal.Range("C16").Value = 0.5

$journal.Range("A17").Value = 45022
$journal.Range("B17").Value = "REST 1 réglement de problème"
$journal.Range("C17").Value = 2.5

# C22 already holds =SUM(C8:C21); it recalculates automatically (20.25 -> 24).

# --- Leave the final selection/scroll position as last touched in the file ---
$planning.Range("G24").Select() | Out-Null
$journal.Range("B19").Select() | Out-Null
